# Auto-generated edit script applying cryptos.xlsx price/volume updates
# and two row re-ordering swaps (rows 44/45 and 49/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.743.82"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").Value = "3.435.22"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'407.02"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'130.18"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +5.43%  "
$ws.Range("D10").Value = "'0.138"
$ws.Range("E10").Value = "  +16.17%  "
$ws.Range("D11").Value = "'41.98"
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "'19.84"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").Value = "3.425.66"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "62.709.88"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "'11.49"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +20.98%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "'84.46"
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("D22").Value = "'315.53"
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("D23").Value = "'12.82"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "'3.17"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "'4.77"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("E28").Value = "  +5.69%  "
$ws.Range("D29").Value = "'2.74"
$ws.Range("E29").Value = "  +7.76%  "
$ws.Range("D30").Value = "'44.47"
$ws.Range("E30").Value = "  +8.78%  "
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "'11.39"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'0.0484"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "'51.86"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +15.62%  "
$ws.Range("D39").Value = "'2.97"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").Value = "'142.75"
$ws.Range("E41").Value = "  +5.60%  "
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").Value = "'16.86"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.91"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "'21.41"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "2.106.81"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'1.97"
$ws.Range("E49").Value = "  +5.48%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").Value = "'1.10"
$ws.Range("E51").Value = "  +31.31%  "
